$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

$ws.Rows("18:18").Insert()

$v = $ws.Range("F21").Validation
$v.Modify(3, 1, 1, "=Intro!`$A`$13:`$A`$19")
Write-Host "modified"
